# "ML model retrained with all data"
# The model's two output/score columns (J = predicted class "r"/col1,
# K = predicted class "s"/col2) are rewritten: every row now reports the
# same retrained-model output (J=0.5, K=1) instead of the previous
# per-row values (J=0.3/row1 "r", K=0.5/row1 "s").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite the whole J and K columns (rows 1-51) with the retrained
# model's constant outputs. Writing plain numbers here also turns the
# former shared-string header cells J1/K1 into numeric cells, matching
# the target workbook (shared strings table ends up empty).
$ws.Range("J1:J51").Value = 0.5
$ws.Range("K1:K51").Value = 1

# Restore the view state Excel recorded after the edit: scrolled down to
# row 19, 90% zoom, and K1:K51 selected with K1 active.
[void]$ws.Range("K1:K51").Select()
$excel.ActiveWindow.Zoom = 90
